# Kuldeep Yadav.xlsx – swap the runs/balls/fours figures recorded for the
# two innings rows (row 2 and row 3) on the "Kuldeep Yadav " sheet.
# Values are stored as text (numbers-as-text), so a leading apostrophe is
# used when assigning through the Range.Value COM property to keep them
# as text instead of letting Excel auto-convert them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: runs/balls/fours -> 12 / 19 / 1
$ws.Range("C2").Value = "'12"
$ws.Range("D2").Value = "'19"
$ws.Range("E2").Value = "'1"

# Row 3: runs/balls/fours -> 1 / 2 / 0
$ws.Range("C3").Value = "'1"
$ws.Range("D3").Value = "'2"
$ws.Range("E3").Value = "'0"
